$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells, styled like the existing header row (e.g. H1).
# Copy H1's formatting (xlPasteAll = -4122) so the new cells reuse the
# same cell style as the rest of the header row, then overwrite the value.
$ws.Range("H1").Copy() | Out-Null

$ws.Range("I1").PasteSpecial(-4122) | Out-Null
$ws.Range("I1").Value = "I0"

$ws.Range("J1").PasteSpecial(-4122) | Out-Null
$ws.Range("J1").Value = "IF"

# New data cells
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 5

$ws.Range("I3").Value = 6
$ws.Range("J3").Value = 6

$ws.Range("I4").Value = 6
$ws.Range("J4").Value = 6
